# Add Splash Screen, MainMenu, and clickable AD Icon
#
# 1) Typography sheet: add a new typography "Iceland_45" (Iceland-Regular.ttf,
#    size 45, bpp 4, fallback "?") as row 7 of the typography table.
# 2) Translation sheet: repurpose the two existing single-use text rows to
#    hold the new "Technology Demonstrator" (splash screen) and "Main Menu"
#    texts (using the new Iceland_45 typography), and add a third row for the
#    new clickable "Analog & Digital" icon text.

$wb = $excel.ActiveWorkbook

# ---- Typography sheet --------------------------------------------------
$wsTypo = $wb.Worksheets.Item("Typography")

$wsTypo.Range("B7").Value = "Iceland_45"
$wsTypo.Range("C7").Value = "Iceland-Regular.ttf"
$wsTypo.Range("D7").Value = 45
$wsTypo.Range("E7").Value = 4
$wsTypo.Range("F7").Value = "?"

# ---- Translation sheet --------------------------------------------------
$wsTrans = $wb.Worksheets.Item("Translation")

# Row 4: was SingleUseId1 / Default / Left / LTR / "Hello IMR"
#  -> SingleUseId5 / Iceland_45 / Left / LTR / "Technology Demonstrator"
$wsTrans.Range("B4").Value = "SingleUseId5"
$wsTrans.Range("C4").Value = "Iceland_45"
$wsTrans.Range("F4").Value = "Technology Demonstrator"

# Row 5: was SingleUseId2 / Default / Left / LTR / "HELLO HAB"
#  -> SingleUseId6 / Iceland_45 / Left / LTR / "Main Menu"
$wsTrans.Range("B5").Value = "SingleUseId6"
$wsTrans.Range("C5").Value = "Iceland_45"
$wsTrans.Range("F5").Value = "Main Menu"

# Row 6 (new): SingleUseId7 / Iceland_45 / Left / LTR / "Analog & Digital"
$wsTrans.Range("B6").Value = "SingleUseId7"
$wsTrans.Range("C6").Value = "Iceland_45"
$wsTrans.Range("D6").Value = "Left"
$wsTrans.Range("E6").Value = "LTR"
$wsTrans.Range("F6").Value = "Analog & Digital"
